$wb = $excel.ActiveWorkbook

# Remove the old "testcase_v1" sheet, keep "testcase_v2" and rename it to "testcase"
$ws1 = $wb.Worksheets.Item("testcase_v1")
$ws1.Delete() | Out-Null
$ws = $wb.Worksheets.Item("testcase_v2")
$ws.Name = "testcase"
$ws.Activate() | Out-Null

# Update row 2 (TC_LOGIN_0001) values
$ws.Range("E2").Value = "phuongtt-auto-stg-01"
$ws.Range("F2").Value = "PhuongTT@12345"
$ws.Range("J2").Value = "phuongtt-auto-stg-01"
$ws.Range("L2").Value = "CHI-LINH-123"

# Update row 3 (TC_LOGIN_0002) values
$ws.Range("B3").Value = "Verify Login Success 02"
$ws.Range("E3").Value = "phuongtt-auto-stg-01"
$ws.Range("F3").Value = "PhuongTT@12345"
$ws.Range("J3").Value = "phuongtt-auto-stg-01"
$ws.Range("L3").Value = "CHI-LINH-123"

# Add hyperlinks on the password cells
$ws.Hyperlinks.Add($ws.Range("F2"), "https://example.com/PhuongTT@12345", "", "", "PhuongTT@12345") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://example.com/PhuongTT@12345", "", "", "PhuongTT@12345") | Out-Null

# Widen the user_name / password columns to fit the new, longer values
$ws.Columns.Item(5).ColumnWidth = 28.94
$ws.Columns.Item(6).ColumnWidth = 30.61

$ws.Range("E25").Select() | Out-Null
